$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last status check" timestamp in the header row
$ws.Range("F1").Value = "Last status check on: 10.02.2022 01:00"

# Row 5 (Makro) price refresh: Cena/Old Cena swapped, Delta Cena and Old Datum
# written out as plain text by the updater script rather than numeric values.
$ws.Range("B5").Value = 35.5
$ws.Range("C5").Value = 35.9
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "-0.4"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = "2022-02-10 01:01:44"
$ws.Range("E5").Style = "Normal"
